$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.931.80'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '3.047.75'
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '582.18'
$ws.Range("E5").Value = '  -1.49%  '
$ws.Range("D6").Value = '150.55'
$ws.Range("E6").Value = '  -2.59%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  -2.20%  '
$ws.Range("D9").Value = '3.051.07'
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("E10").Value = '  -3.06%  '
$ws.Range("D11").Value = '5.79'
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("E12").Value = '  -1.78%  '
$ws.Range("D13").Value = '0.0000232'
$ws.Range("E13").Value = '  -3.31%  '
$ws.Range("D14").Value = '35.78'
$ws.Range("E14").Value = '  -4.31%  '
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").Value = '3.553.94'
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("D17").Value = '7.09'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = '62.981.30'
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("D19").Value = '3.048.54'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("D20").Value = '477.37'
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").Value = '14.20'
$ws.Range("E21").Value = '  -2.64%  '
$ws.Range("D22").Value = '0.702'
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").Value = '7.49'
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("E24").Value = '  -1.52%  '
$ws.Range("D25").Value = '81.38'
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("D26").Value = '12.56'
$ws.Range("E26").Value = '  -2.45%  '
$ws.Range("E27").Value = '  +5.13%  '
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").Value = '7.25'
$ws.Range("E30").Value = '  -1.36%  '
$ws.Range("E31").Value = '  -1.56%  '
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("D33").Value = '27.65'
$ws.Range("E33").Value = '  +1.92%  '
$ws.Range("E34").Value = '  -3.70%  '
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("D36").Value = '0.0₃0804'
$ws.Range("E36").Value = '  -5.03%  '
$ws.Range("E37").Value = '  -3.31%  '
$ws.Range("E38").Value = '  -1.72%  '
$ws.Range("D39").Value = '3.05'
$ws.Range("E39").Value = '  -9.59%  '
$ws.Range("D40").Value = '50.18'
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("D41").Value = '9.09'
$ws.Range("E41").Value = '  -1.84%  '
$ws.Range("D42").Value = '425.19'
$ws.Range("E42").Value = '  -4.03%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  +2.34%  '
$ws.Range("D45").Value = '2.829.15'
$ws.Range("E45").Value = '  +1.05%  '
$ws.Range("D46").Value = '0.0359'
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("D47").Value = '38.00'
$ws.Range("E47").Value = '  -5.11%  '
$ws.Range("D48").Value = '127.36'
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("D49").Value = '0.999'
$ws.Range("D50").Value = '25.06'
$ws.Range("E50").Value = '  -2.34%  '
$ws.Range("E51").Value = '  -0.77%  '
